$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The product-name column (A2:A21) used to be a single merged cell holding
# "CUSTOM AL QURAN MENGENANG/WAFAT 40/100/1000 HARI". The edit un-merges it
# and repeats the same label down every row instead, carrying over A2's
# look (minus the forced horizontal centering that merged cells get).
# ---------------------------------------------------------------------------

$label = $ws.Range("A2").Value2

# Un-merge A2:A21 (also drops the <mergeCells> entry).
$ws.Range("A2:A21").UnMerge()

# Fill every row of column A with the same product label.
for ($r = 3; $r -le 21; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $label
}

# A standalone (non-merged) cell doesn't need horizontal centering anymore;
# keep it vertically centered only.
$ws.Range("A2").HorizontalAlignment = 1   # xlGeneral

# Copy A2's resulting format (border + alignment) down onto A3:A21 so every
# row in the column matches exactly.
$ws.Range("A2").Copy()
$ws.Range("A3:A21").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Minor view-state touch-ups captured in the saved file.
# ---------------------------------------------------------------------------
$ws.Range("B22").Select() | Out-Null
